$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the last data row (row 67): remove the ID/meter-number/meter-type
# values entirely and blank out the two remaining (styled) cells, mirroring
# a "get data to model" re-pull that dropped the last reading.
$ws.Range("A67:E67").ClearContents()

# Scroll the view down so row 55 is at the top and leave the last row
# selected/active, matching the author's on-screen position when saving.
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A67:G67").Select()
